# Applies the two substantive content changes described by the commit:
#  1. Fix a typo: "a bit of work clean up work" -> "a bit of cleanup work",
#     and add a comma: "reset service values so" -> "reset service values, so".
#  2. Add a new paragraph about the Student record being saved/updated by
#     the AssessmentEngine, right after the paragraph that ends with
#     "...decides which path to take next in the algorithm."
#
# (The rest of the upstream diff is just Word's proofing-error markers
# <w:proofErr .../> being added/removed/re-merged around unchanged runs --
# invisible to the rendered text, so no Find/Replace is needed for those.)

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1 (Replace:=wdReplaceOne)
$wdFindContinue = 1
$wdReplaceOne = 1

# --- Change 1a: "a bit of work clean up work" -> "a bit of cleanup work" ---
$d.Content.Find.Execute(
    "there is a bit of work clean up work that needs to be done",
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "there is a bit of cleanup work that needs to be done",
    $wdReplaceOne) | Out-Null

# --- Change 1b: "reset service values so" -> "reset service values, so" ---
$d.Content.Find.Execute(
    "to reset service values so the service is ready",
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "to reset service values, so the service is ready",
    $wdReplaceOne) | Out-Null

# --- Change 2: insert the new "Student record" paragraph ---
$oldTail = "awaits this result from the server before it decides which path to take next in the algorithm."
$newParagraphText = "Finally, the code for saving and updating the Student record is in the AssessmentEngine. These http requests are kept here because the Student record is first created and saved when the assessment begins. Additionally, any updates that happen to a student record is done so after searching for TakenAssessment records, which is also managed by the AssessmentEngine."
$newTail = $oldTail + "`r" + $newParagraphText

$d.Content.Find.Execute(
    $oldTail,
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    $newTail,
    $wdReplaceOne) | Out-Null

Write-Output "Edits applied."
